$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 206, shifting existing rows 206+ down by 2.
$ws.Rows("206:207").Insert()

$ws.Range("A206").Value = 9
$ws.Range("B206").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C206").Value = "Metropolitana"
$ws.Range("D206").Value = 44529
$ws.Range("E206").Value = 13
$ws.Range("F206").Value = 100114014
$ws.Range("G206").Value = "Betarraga"
$ws.Range("H206").Value = "Sin especificar"
$ws.Range("I206").Value = "Primera"
$ws.Range("J206").Value = 3400
$ws.Range("K206").Value = 100
$ws.Range("L206").Value = 110
$ws.Range("M206").Value = 105
$ws.Range("N206").Value = "`$/unidad"
$ws.Range("O206").Value = "Región Metropolitana"
$ws.Range("P206").Value = 105
$ws.Range("Q206").Value = 1
$ws.Range("R206").Value = "Hortaliza"

$ws.Range("A207").Value = 9
$ws.Range("B207").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C207").Value = "Metropolitana"
$ws.Range("D207").Value = 44529
$ws.Range("E207").Value = 13
$ws.Range("F207").Value = 100114014
$ws.Range("G207").Value = "Betarraga"
$ws.Range("H207").Value = "Sin especificar"
$ws.Range("I207").Value = "Segunda"
$ws.Range("J207").Value = 1600
$ws.Range("K207").Value = 80
$ws.Range("L207").Value = 90
$ws.Range("M207").Value = 85
$ws.Range("N207").Value = "`$/unidad"
$ws.Range("O207").Value = "Región Metropolitana"
$ws.Range("P207").Value = 85
$ws.Range("Q207").Value = 1
$ws.Range("R207").Value = "Hortaliza"
